$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task list was reordered: the item that was in row 8
# ("Error handling strategy in Model Compiler", estimate 7) moved up to
# become the new row 2, and the items that were previously in rows 2-7
# each shifted down by one row. Rows 9 and 10 are unchanged.

$tasks = @(
    "Error handling strategy in Model Compiler",
    "Refactoring - we need consistency across the board - All 3 apps",
    "Make exporter a GUP.  Build UI & hook data into the max files",
    "Make the path from Max->Model Viewer seamless",
    "Add full screen support",
    "Add camera controls to model viewer",
    "Textured surfaces",
    "Complete the Rorn Maths library",
    "Revise, understand and document the view and projection matrix builds"
)

$estimates = @(7, 21, 14, 4, 3, 5, 35, 35, 7)

for ($i = 0; $i -lt $tasks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tasks[$i]
    $ws.Cells.Item($row, 2).Value = $estimates[$i]
}
